# Refresh the Leve profit-calculation columns (H:N) on the Bahamut_Profits
# crafting-class sheets, mirroring the scheduled market-data runner.
# currentAveragePrice(NQ/HQ) + LevePrice(NQ/HQ) + LeveProfit(NQ/HQ) move to
# their latest observed values; a couple of rows also gain/lose a cell
# because their LeveProfit is blank (N/A) under the new figures.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 2000.05
$ws.Range("I40").Value = 1001
$ws.Range("J40").Value = 2052.6316
$ws.Range("K40").Value = 1001
$ws.Range("L40").Value = 2052.6316
$ws.Range("M40").Value = -826
$ws.Range("N40").Value = -2402.6316

# Row 51
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").ClearContents()

# Row 62
$ws.Range("H62").Value = 54997.332
$ws.Range("I62").Value = 105474.4
$ws.Range("J62").Value = 9109.091
$ws.Range("K62").Value = 105474.4
$ws.Range("L62").Value = 9109.091
$ws.Range("M62").Value = -104850.4
$ws.Range("N62").Value = -10357.091

# Row 64
$ws.Range("H64").Value = 6714.2856
$ws.Range("I64").Value = 9050
$ws.Range("J64").Value = 3600
$ws.Range("K64").Value = 9050
$ws.Range("L64").Value = 3600
$ws.Range("M64").Value = -8802
$ws.Range("N64").Value = -4096

# Row 65
$ws.Range("H65").Value = 54997.332
$ws.Range("I65").Value = 105474.4
$ws.Range("J65").Value = 9109.091
$ws.Range("K65").Value = 527372
$ws.Range("L65").Value = 45545.455
$ws.Range("M65").Value = -524252
$ws.Range("N65").Value = -51785.455

# Row 67
$ws.Range("H67").Value = 6714.2856
$ws.Range("I67").Value = 9050
$ws.Range("J67").Value = 3600
$ws.Range("K67").Value = 9050
$ws.Range("L67").Value = 3600
$ws.Range("M67").Value = -8192
$ws.Range("N67").Value = -5316

# Row 76
$ws.Range("H76").Value = 78950110
$ws.Range("I76").Value = 88238010
$ws.Range("J76").Value = 3000
$ws.Range("K76").Value = 88238010
$ws.Range("L76").Value = 3000
$ws.Range("M76").Value = -88237695
$ws.Range("N76").Value = -3630

# Row 79
$ws.Range("H79").Value = 78950110
$ws.Range("I79").Value = 88238010
$ws.Range("J79").Value = 3000
$ws.Range("K79").Value = 88238010
$ws.Range("L79").Value = 3000
$ws.Range("M79").Value = -88236918
$ws.Range("N79").Value = -5184

# Row 112
$ws.Range("H112").Value = 3997.3684
$ws.Range("J112").Value = 4151.4287
$ws.Range("L112").Value = 12454.2861
$ws.Range("N112").Value = -14670.2861

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1202.871
$ws.Range("I45").Value = 1153.125
$ws.Range("J45").Value = 1373.4286
$ws.Range("K45").Value = 1153.125
$ws.Range("L45").Value = 1373.4286
$ws.Range("M45").Value = -776.125
$ws.Range("N45").Value = -2127.4286

# Row 63
$ws.Range("H63").Value = 2715.6667
$ws.Range("I63").Value = 2145.1428
$ws.Range("K63").Value = 2145.1428
$ws.Range("M63").Value = -1459.1428

# Row 66
$ws.Range("H66").Value = 2715.6667
$ws.Range("I66").Value = 2145.1428
$ws.Range("K66").Value = 10725.714
$ws.Range("M66").Value = -7293.714

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 2940
$ws.Range("I105").Value = 3064.3
$ws.Range("J105").Value = 2473.875
$ws.Range("K105").Value = 3064.3
$ws.Range("L105").Value = 2473.875
$ws.Range("M105").Value = -1317.3
$ws.Range("N105").Value = -5967.875

# Row 107
$ws.Range("H107").Value = 7561.45
$ws.Range("I107").Value = 1103
$ws.Range("J107").Value = 19555.715
$ws.Range("K107").Value = 1103
$ws.Range("L107").Value = 19555.715
$ws.Range("M107").Value = 817
$ws.Range("N107").Value = -23395.715

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 1170.25
$ws.Range("I22").Value = 1591.5714
$ws.Range("K22").Value = 1591.5714
$ws.Range("M22").Value = -1241.5714

# Row 62
$ws.Range("H62").Value = 5512.2144
$ws.Range("I62").Value = 5931.6665
$ws.Range("J62").Value = 4757.2
$ws.Range("K62").Value = 5931.6665
$ws.Range("L62").Value = 4757.2
$ws.Range("M62").Value = -5307.6665
$ws.Range("N62").Value = -6005.2

# Row 65
$ws.Range("H65").Value = 5512.2144
$ws.Range("I65").Value = 5931.6665
$ws.Range("J65").Value = 4757.2
$ws.Range("K65").Value = 29658.3325
$ws.Range("L65").Value = 23786
$ws.Range("M65").Value = -26538.3325
$ws.Range("N65").Value = -30026

# Row 105
$ws.Range("H105").Value = 949.8333
$ws.Range("I105").Value = 949.8333
$ws.Range("K105").Value = 949.8333
$ws.Range("M105").Value = 797.1667

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 129
$ws.Range("H129").Value = 51461.19
$ws.Range("I129").Value = 1736.6666
$ws.Range("J129").Value = 71351
$ws.Range("K129").Value = 5209.9998
$ws.Range("L129").Value = 214053
$ws.Range("M129").Value = -209.9997999999996
$ws.Range("N129").Value = -224053

# Row 131
$ws.Range("H131").Value = 135755.4
$ws.Range("J131").Value = 85460.914
$ws.Range("L131").Value = 256382.742
$ws.Range("N131").Value = -266462.742

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 14000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 14000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 14000
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -14302

# Row 70
$ws.Range("H70").Value = 16949.967
$ws.Range("I70").Value = 30399.334
$ws.Range("J70").Value = 4341.1875
$ws.Range("K70").Value = 30399.334
$ws.Range("L70").Value = 4341.1875
$ws.Range("M70").Value = -30129.334
$ws.Range("N70").Value = -4881.1875

# Row 73
$ws.Range("H73").Value = 16949.967
$ws.Range("I73").Value = 30399.334
$ws.Range("J73").Value = 4341.1875
$ws.Range("K73").Value = 30399.334
$ws.Range("L73").Value = 4341.1875
$ws.Range("M73").Value = -29463.334
$ws.Range("N73").Value = -6213.1875

# Row 80
$ws.Range("H80").Value = 3650.4546
$ws.Range("I80").Value = 4451
$ws.Range("J80").Value = 2983.3333
$ws.Range("K80").Value = 4451
$ws.Range("L80").Value = 2983.3333
$ws.Range("M80").Value = -3453
$ws.Range("N80").Value = -4979.3333

# Row 83
$ws.Range("H83").Value = 3650.4546
$ws.Range("I83").Value = 4451
$ws.Range("J83").Value = 2983.3333
$ws.Range("K83").Value = 22255
$ws.Range("L83").Value = 14916.6665
$ws.Range("M83").Value = -17263
$ws.Range("N83").Value = -24900.6665

# Row 113
$ws.Range("H113").Value = 3309.9614
$ws.Range("I113").Value = 3004.762
$ws.Range("J113").Value = 4591.8
$ws.Range("K113").Value = 3004.762
$ws.Range("L113").Value = 4591.8
$ws.Range("M113").Value = -834.7620000000002
$ws.Range("N113").Value = -8931.799999999999

# Row 122
$ws.Range("H122").Value = 732733.75
$ws.Range("I122").Value = 1197558.2
$ws.Range("J122").Value = 2295.1428
$ws.Range("K122").Value = 3592674.6
$ws.Range("L122").Value = 6885.428400000001
$ws.Range("M122").Value = -3590224.6
$ws.Range("N122").Value = -11785.4284

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 6404.2085
$ws.Range("I136").Value = 3550.125
$ws.Range("K136").Value = 10650.375
$ws.Range("M136").Value = -8100.375
